$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new inventory row (id 40) for hEGF
$ws.Range("A41").Value = 40
$ws.Range("B41").Value = "hEGF"
$ws.Range("C41").Value = "E9644-.5MG"
$ws.Range("D41").Value = "https://www.sigmaaldrich.com/MX/es/product/sigma/e9644?utm_source=google&utm_medium=cpc&utm_campaign=22179178721&utm_content=177793360441&gad_source=1&gad_campaignid=22179178721&gbraid=0AAAAAD8kLQQijHqI3ftWe1h-KzX0veMTN&gclid=CjwKCAjwtfvEBhAmEiwA-DsKjncC52MuEisnbBttGKNLUW_jnD-jpKIWAwJoS3ZTUteNgQdjhKakFxoCKgUQAvD_BwE"
$ws.Range("E41").Value = "Small Molecule"
$ws.Range("F41").Value = "Aliquot"
$ws.Range("G41").Value = -30
$ws.Range("H41").Value = 9
$ws.Range("I41").Value = "Opened"
$ws.Range("J41").Value = 52
$ws.Range("K41").Value = "1 mL"

# Update the sheet's view to reflect the scrolled/selected state after the edit
$ws.Activate()
[void]$ws.Range("M45").Select()
